$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9
$ws.Range("A9").Value = "PCEDGC96"
$ws.Range("B9").Value = "rtdsm"
$ws.Range("C9").Value = "RCOND"
$ws.Range("D9").Value = "rcondMvQd.xlsx"

# Row 10 (filename filled before the variable name, matching original authoring order)
$ws.Range("A10").Value = "PCENDC96"
$ws.Range("B10").Value = "rtdsm"
$ws.Range("D10").Value = "rconndMvQd.xlsx"
$ws.Range("C10").Value = "RCONND"

# Row 11 (filename filled before the variable name, matching original authoring order)
$ws.Range("A11").Value = "PCESC96"
$ws.Range("B11").Value = "rtdsm"
$ws.Range("D11").Value = "rconsMvQd.xlsx"
$ws.Range("C11").Value = "RCONS"

$ws.Range("C15").Select()
